# The document's logos (header + two footers) have their inline-picture
# "name" attribute (wp:docPr/@name and pic:cNvPr/@name) out of sync with
# the embedded media filenames used elsewhere in the package. Relabel
# them so the Pearson logos (currently "image2.png") become "image1.png",
# and the BTEC logo (currently "image1.jpg") becomes "image2.jpg".

$d = $word.ActiveDocument
$sec = $d.Sections.First

# Header: BTEC logo -> rename image1.jpg to image2.jpg
$hdr = $sec.Headers.Item(2)
$hdr.Range.InlineShapes.Item(1).Name = "image2.jpg"

# Footer (page 1 style): Pearson logo -> rename image2.png to image1.png
$ftr1 = $sec.Footers.Item(1)
$ftr1.Range.InlineShapes.Item(1).Name = "image1.png"

# Footer (other pages style): Pearson logo -> rename image2.png to image1.png
$ftr2 = $sec.Footers.Item(2)
$ftr2.Range.InlineShapes.Item(1).Name = "image1.png"
